$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New section header row (row 12), mirrors the row 1 header ---
$ws.Range("A12:G12").Merge() | Out-Null
$ws.Range("A1:G1").Copy()
$ws.Range("A12:G12").PasteSpecial(-4122)
$ws.Range("A12").Value = "SwipeOutController"

# --- Copy the column header row (row 2) down to row 13 ---
$ws.Range("A2:G2").Copy()
$ws.Range("A13:G13").PasteSpecial(-4122)
$ws.Range("A2:G2").Copy()
$ws.Range("A13:G13").PasteSpecial(-4163)

# --- Copy the body rows (3-10) down to rows 14-21 (values + formats) ---
$ws.Range("A3:G10").Copy()
$ws.Range("A14:G21").PasteSpecial(-4122)
$ws.Range("A3:G10").Copy()
$ws.Range("A14:G21").PasteSpecial(-4163)

# Match source row heights for the wrapped-text rows that were copied
$ws.Rows("15").RowHeight = $ws.Rows("4").RowHeight
$ws.Rows("19").RowHeight = $ws.Rows("8").RowHeight
$ws.Rows("20").RowHeight = $ws.Rows("9").RowHeight

# Row 19 (step 6) gets new "checkout" wording instead of the swipe-in text
$ws.Range("B19").Value = "User allowed to get out:" + [char]10 + "When Balance meets minimum account balance criteria"

# --- New row 22: additional test step; only A22:B22 get the new border/fill style ---
$ws.Range("A21:G21").Copy()
$ws.Range("A22:G22").PasteSpecial(-4122)
$ws.Range("A22:B22").Interior.ColorIndex = -4142
$ws.Range("B22").WrapText = $true
$ws.Range("A22").Value = 9
$ws.Range("B22").Value = "Actual Fare deducted on the basis of days and Number of Stations"
$ws.Rows("22").RowHeight = 22.5

# --- Update view: scroll/selection to the newly added scenario block ---
$ws.Range("A13:G22").Select() | Out-Null

$excel.CutCopyMode = $false
